# Update TrialsSetup 2026-02-17 16:00
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Progress value for row 3: 12.5 -> 25
$ws.Range("C3").Value = 25

# Days remaining value for row 6: 3 -> 2
$ws.Range("B6").Value = 2

# Days remaining value for row 8: 23 -> 22
$ws.Range("B8").Value = 22

$wb.Save()
